$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Evolving Cyber Security Landscape in a Digital Age") | Out-Null
$rng.Text = "Unmasking the Realm of Matter: A Journey into Chemistry's Enchanting Laboratory"

# ---------------------------------------------------------------------------
# 2) Author name
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Dylan Jones") | Out-Null
$rng.Text = "Emily Richards"

# ---------------------------------------------------------------------------
# 3) Email address (two runs: the account part and the TLD part)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("dylanjones@emailworld") | Out-Null
$rng.Text = "emilyrichards@highschool"

$rng = $d.Content
$rng.Find.Execute("dylanjones@emailworld") | Out-Null
$rng.Collapse(0) | Out-Null
$rng.MoveEnd(1, 1) | Out-Null
$rng.MoveStart(1, 1) | Out-Null
$rng.Text = "edu"

# ---------------------------------------------------------------------------
# 4) Body paragraph - three "blocks" separated by manual line breaks
#    (double <w:br/>). Each block is replaced in full (original sentences +
#    new sentences appended by the edit) so the line-break landmarks stay put.
# ---------------------------------------------------------------------------

# --- Block 1 ---
$startRng = $d.Content
$startRng.Find.Execute("In a world driven by interconnected technologies") | Out-Null
$startPos = $startRng.Start

$endRng = $d.Content
$endRng.Find.Execute("This essay delves into the intricacies of the ever-changing cyber security landscape, exploring the challenges and advancements that shape this rapidly evolving field") | Out-Null
$endPos = $endRng.End

$blockRng = $d.Range($startPos, $endPos)
$blockRng.Text = "The realm of chemistry is akin to a vast and mystical kingdom, holding within its depths a boundless tapestry of intricate interactions and transformative processes. Like an arcane sorcerer, the chemist wields the power of elements and molecules, unraveling the enigmatic symphony of nature's very fabric. Through the prism of chemical reactions, we witness the symphony of atoms dancing in harmonious motion, revealing the profound interconnectedness of all matter that shapes our universe. In our laboratory crucible, we embark on an alchemical quest to decode the enigmatic secrets hidden within the molecular realm, unveiling the fundamental principles that govern the intricate play of elements. As we delve further into this wondrous realm, we unravel the mysteries of chemical reactions, unmasking the underlying order amidst apparent chaos, and illuminating the deep resonance between chemistry and our everyday world"

# --- Block 2 ---
$startRng = $d.Content
$startRng.Find.Execute("The increasing reliance on technology and the proliferation of digital data") | Out-Null
$startPos = $startRng.Start

$endRng = $d.Content
$endRng.Find.Execute("These attacks can result in data breaches, financial fraud, disruption of critical services, and reputational damage") | Out-Null
$endPos = $endRng.End

$blockRng = $d.Range($startPos, $endPos)
$blockRng.Text = "With each experiment conducted, like a detective meticulously solving a puzzle, we piece together the chemistry's enchanting story. We explore the remarkable properties of elements, venturing into the vibrant world of periodic trends and unraveling the mysteries of atomic structure. We unravel the profound implications of chemical bonding, revealing the delicate balance between attraction and repulsion that shapes molecular architecture and governs the reactivity of substances. Through our investigations, we uncover the dynamic realm of chemical reactions, witnessing the transformation of reactants into products, and glimpse the energetic interplay that drives these processes"

# --- Block 3 ---
$startRng = $d.Content
$startRng.Find.Execute("Recognizing the gravity of these threats") | Out-Null
$startPos = $startRng.Start

$endRng = $d.Content
$endRng.Find.Execute("The development of security frameworks and standards has also contributed to a more structured and systematic approach to cyber security management") | Out-Null
$endPos = $endRng.End

$blockRng = $d.Range($startPos, $endPos)
$blockRng.Text = "In unraveling the intricacies of chemical reactions, we discover the profound influence of temperature, pressure, and catalysis, orchestrating the dance of molecules and dictating the pace and direction of change. We venture into the fascinating realm of equilibrium, where opposing forces clash in delicate balance, revealing the interplay of spontaneity and stability. Furthermore, we elucidate the concepts of acids and bases, unveiling their pivotal role in countless natural and industrial processes, shaping the world around us in myriad ways"

# ---------------------------------------------------------------------------
# 5) "Summary" heading gains a lastRenderedPageBreak marker ahead of its text
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Summary") | Out-Null
$rng.Select() | Out-Null
$word.Selection.Range.InsertParagraphBefore() | Out-Null

# ---------------------------------------------------------------------------
# 6) Summary body paragraph
# ---------------------------------------------------------------------------
$startRng = $d.Content
$startRng.Find.Execute("The cyber security landscape is a dynamic and ever-evolving domain") | Out-Null
$startPos = $startRng.Start

$endRng = $d.Content
$endRng.Find.Execute("By embracing innovation and adopting best practices, we can effectively navigate the evolving cyber security landscape and safeguard our digital world") | Out-Null
$endPos = $endRng.End

$blockRng = $d.Range($startPos, $endPos)
$blockRng.Text = "In this exploration of chemistry's enchanting realm, we have ventured into the depths of matter's innermost sanctum, unmasking the enigmatic symphony of chemical reactions and unveiling the underlying principles that govern the intricate play of elements and molecules. Through the crucible of laboratory experiments, we have pieced together the chemistry's captivating story, unraveling the mysteries of elements, bonding, and reactions. We have witnessed the transformative power of chemistry, its profound influence on our world, and its limitless potential to unlock the secrets of nature. As we continue our journey into this realm of molecular enchantment, we stand on the precipice of even greater discoveries, ready to unravel the mysteries that still lie hidden, awaiting our exploration"

# ---------------------------------------------------------------------------
# 7) New empty paragraph at the very end of the document body
# ---------------------------------------------------------------------------
$endRng = $d.Content
$endRng.Collapse(0) | Out-Null
$endRng.InsertParagraphAfter() | Out-Null
